$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three new blank columns before column B. This shifts the
#    existing B:E data (and any per-cell formatting) right to E:H,
#    exactly like the target workbook.
# ---------------------------------------------------------------------
$ws.Columns("B:D").Insert()

# ---------------------------------------------------------------------
# 2. New header row values for the freshly inserted columns.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# ---------------------------------------------------------------------
# 3. Fill the new columns with "UN" for every existing data row
#    (rows 2-27), matching the rest of the table.
# ---------------------------------------------------------------------
For ($r = 2; $r -le 27; $r++) {
    $ws.Range("B" + $r).Value = "UN"
    $ws.Range("C" + $r).Value = "UN"
    $ws.Range("D" + $r).Value = "UN"
}

# ---------------------------------------------------------------------
# 4. Row 11 (ValuEngine) gets the new rating-change notes in C11/D11,
#    highlighted the same way the existing upgrade note in E11 is.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = "6/22/2018,Upgrades,Sell -> Hold,"
$ws.Range("D11").Value = "6/22/2018,Upgrades,Sell -> Hold,"
$ws.Range("C11").Interior.Color = 13434828
$ws.Range("D11").Interior.Color = 13434828

# ---------------------------------------------------------------------
# 5. Append the two new analyst rows at the bottom of the table.
# ---------------------------------------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
